$wb = $excel.ActiveWorkbook

# --- Sheet "Users": contact person changed ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "James Craven"
[void]$wsUsers.Range("E3").Select()

# --- Sheet "ActivityColumns": new "." / ".." rows inserted, "Description" column removed ---
$wsCols = $wb.Worksheets.Item("ActivityColumns")
$wsCols.Range("A3").Value = "."
$wsCols.Range("A4").Value = ".."
$wsCols.Range("A5").Value = "Date"
$wsCols.Range("A6").Value = "Company Name"
$wsCols.Range("A7").Value = "Type"
$wsCols.Range("A8").Value = "Tier"
$wsCols.Range("A9").Value = "Event/Task Type"
$wsCols.Range("A10").Value = "HL Contact"
$wsCols.Range("A11").Value = "Subject"
$wsCols.Range("A12").Value = "Meeting/Call Notes"
$wsCols.Range("A13").Value = "External Contact"

# ActivityColumns becomes the active tab/selection
[void]$wsCols.Range("A5").Select()
